$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Active")

# Insert two new rows at row 2, pushing existing data down
$ws.Rows.Item(2).Resize(2).Insert()

# Populate the two newly inserted rows
$ws.Cells.Item(2, 2).Value = "Supergirl - Anna Naklab ft. Allen Farben"
$ws.Cells.Item(2, 3).Value = "https://www.youtube.com/watch?v=swBR4QnO3yE"

$ws.Cells.Item(3, 2).Value = "Maroon 5, Wiz Khalifa – Payphone"
$ws.Cells.Item(3, 3).Value = "https://www.youtube.com/watch?v=bbdsIR4UHDg"
